$d = $word.ActiveDocument

# Update the header date
$d.Content.Find.Execute("2026-02-24 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-02-25 Wednesday", 2)

# Helper to replace the content of a table cell while preserving formatting
function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    $rng.End = $rng.End - 1
    $rng.Text = $newText
}

$t = $d.Tables.Item(1)

# Row 1 (table row index 1)
Set-CellText $t 1 1 "79×26="
Set-CellText $t 1 2 "43×92="
Set-CellText $t 1 3 "47×90="
Set-CellText $t 1 4 "21×30="
Set-CellText $t 1 5 "62×13="

# Row 2 (table row index 5)
Set-CellText $t 5 1 "61×54="
Set-CellText $t 5 2 "13×36="
Set-CellText $t 5 3 "54×76="
Set-CellText $t 5 4 "68×41="
Set-CellText $t 5 5 "34×62="

# Row 3 (table row index 10)
Set-CellText $t 10 1 "15×66="
Set-CellText $t 10 2 "41×53="
Set-CellText $t 10 3 "85×42="
Set-CellText $t 10 4 "71×74="
Set-CellText $t 10 5 "38×68="

# Row 4 (table row index 15)
Set-CellText $t 15 1 "95×52="
Set-CellText $t 15 2 "53×90="
Set-CellText $t 15 3 "91×74="
Set-CellText $t 15 4 "40×88="
Set-CellText $t 15 5 "38×49="

# Row 5 (table row index 20)
Set-CellText $t 20 1 "45×78="
Set-CellText $t 20 2 "75×80="
Set-CellText $t 20 3 "48×28="
Set-CellText $t 20 4 "45×18="
Set-CellText $t 20 5 "72×81="
